# Apply crypto price/volume updates scraped on Fri Apr  5 12:24:37 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '66.532.75'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.241.08'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'577.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').Value = "'170.50"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.59%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = "'0.574"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').Value = '3.235.99'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('E10').Value = '  -7.13%  '
$ws.Range('D11').Value = "'0.566"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.86%  '
$ws.Range('D12').Value = "'44.42"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.69%  '
$ws.Range('D13').Value = "'0.0000267"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = "'676.57"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.39%  '
$ws.Range('D15').Value = '3.766.73'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').Value = "'8.17"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('D17').Value = '66.611.38'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '3.234.27'
$ws.Range('E19').Value = '  -3.16%  '
$ws.Range('D20').Value = "'17.04"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.86%  '
$ws.Range('D21').Value = "'10.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.19%  '
$ws.Range('D22').Value = "'0.868"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').Value = "'16.82"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.25%  '
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').Value = "'96.22"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.09%  '
$ws.Range('E26').Value = '  -4.74%  '
$ws.Range('D27').Value = "'2.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.71%  '
$ws.Range('D28').Value = "'8.91"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.85%  '
$ws.Range('D29').Value = "'32.16"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = "'8.19"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.20%  '
$ws.Range('D31').Value = "'6.62"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.72%  '
$ws.Range('D32').Value = "'566.38"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.20%  '
$ws.Range('D33').Value = "'10.75"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.34%  '
$ws.Range('D34').Value = '3.790.92'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('E36').Value = '  -4.43%  '
$ws.Range('D37').Value = "'54.99"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.97%  '
$ws.Range('D38').Value = "'3.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -17.89%  '
$ws.Range('D39').Value = "'0.128"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').Value = "'2.56"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.38%  '
$ws.Range('D41').Value = "'31.18"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.66%  '
$ws.Range('D42').Value = "'3.28"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.63%  '
$ws.Range('D43').Value = '0.0₃0645'
$ws.Range('E43').Value = '  -8.18%  '
$ws.Range('D44').Value = "'2.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.89%  '
$ws.Range('D45').Value = "'0.320"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.99%  '
$ws.Range('D46').Value = "'0.0394"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.71%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('E49').Value = '  -2.34%  '
$ws.Range('D50').Value = "'1.30"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.92%  '
$ws.Range('D51').Value = "'127.09"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.91%  '
